# Add a "bio"/unit picklist block (name/code/prices/unit/min-qty/stock) with
# input validation (dropdown list, non-negative whole numbers) and highlight
# the required header columns, per "add charactor limit about bio input".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): unlock every header cell, and additionally
#    highlight the *required* columns (name, prices, unit, min qty, stock)
#    with the light-yellow fill already used elsewhere in the workbook.
# ---------------------------------------------------------------------------
$requiredHeaderCols = @("A1","D1","E1","F1","G1","H1","M1","N1","P1")
$optionalHeaderCols = @("B1","C1","I1","J1","K1","L1","O1")

foreach ($addr in $requiredHeaderCols) {
    $c = $ws.Range($addr)
    $c.Interior.Color = 13434879   # RGB(255,255,204) - existing light-yellow fill
    $c.Locked = $false
}
foreach ($addr in $optionalHeaderCols) {
    $ws.Range($addr).Locked = $false
}

# ---------------------------------------------------------------------------
# 2) Sample data rows 2-5.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "测试"
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = "xxxx"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 0.01
$ws.Range("J2").Value = 1.01
$ws.Range("K2").Value = 2.01
$ws.Range("L2").Value = 3.01
$ws.Range("M2").Value = "Unidad"
$ws.Range("N2").Value = 1
$ws.Range("P2").Value = 11111

$ws.Range("A3").Value = "测试2"
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = "xxxx"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 8
$ws.Range("H3").Value = 11
$ws.Range("I3").Value = 1.01
$ws.Range("J3").Value = 2.01
$ws.Range("K3").Value = 3.01
$ws.Range("L3").Value = 4.01
$ws.Range("M3").Value = "Caja"
$ws.Range("N3").Value = 12
$ws.Range("P3").Value = 22222

$ws.Range("A4").Value = "测试3"
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = "xxxx"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 2.01
$ws.Range("J4").Value = 3.01
$ws.Range("K4").Value = 4.01
$ws.Range("L4").Value = 5.01
$ws.Range("M4").Value = "Paquete"
$ws.Range("N4").Value = 100
$ws.Range("P4").Value = 33333

$ws.Range("A5").Value = "test"
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = "xxxxx"
$ws.Range("E5").Value = 55
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 3.01
$ws.Range("J5").Value = 4.01
$ws.Range("K5").Value = 5.01
$ws.Range("L5").Value = 6.01
$ws.Range("M5").Value = "Caja"
$ws.Range("N5").Value = 1000
$ws.Range("P5").Value = 44444

# Hidden helper list (column S) backing the unit dropdown.
$ws.Range("S2").Value = "Unidad"
$ws.Range("S3").Value = "Caja"
$ws.Range("S4").Value = "Paquete"
$ws.Range("S5").Value = "Palet"

# ---------------------------------------------------------------------------
# 3) Number formats for the new numeric columns.
# ---------------------------------------------------------------------------
$ws.Range("E2:H5").NumberFormat = "0.000_ "
$ws.Range("I2:L5").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# 4) Unlock all of the newly entered data cells so users can fill the sheet
#    in even when the workbook is protected.
# ---------------------------------------------------------------------------
$ws.Range("A2:A5").Locked = $false
$ws.Range("B2:B5").Locked = $false
$ws.Range("D2:D5").Locked = $false
$ws.Range("E2:H5").Locked = $false
$ws.Range("I2:L5").Locked = $false
$ws.Range("M2:M5").Locked = $false
$ws.Range("N2:N5").Locked = $false
$ws.Range("P2:P5").Locked = $false

# Helper list + the stray formatted cell stay locked, but keep protection
# formatting applied (mirrors the source workbook).
$ws.Range("S2:S5").Locked = $true
$ws.Range("O15").Locked = $true

# ---------------------------------------------------------------------------
# 5) Data validation.
# ---------------------------------------------------------------------------
# M1: informational prompt only, no restriction.
$v1 = $ws.Range("M1").Validation
$v1.Add(0, 1, 1)
$v1.InputTitle = "Solo se puede elejir uno"
$v1.InputMessage = "Unidad, Caja, Paquete, Palet"
$v1.ShowInput = $true
$v1.ShowError = $true

# M2:M1048576: dropdown list sourced from the hidden helper column S.
$v2 = $ws.Range("M2:M1048576").Validation
$v2.Add(3, 1, 1, '=$S$2:$S$5')
$v2.InputTitle = "Solo se puede elejir uno"
$v2.InputMessage = "Unidad, Caja, Paquete, Palet"
$v2.ShowInput = $true
$v2.ShowError = $true

# N2:N1048576: whole numbers >= 0.
$v3 = $ws.Range("N2:N1048576").Validation
$v3.Add(1, 1, 7, 0)
$v3.ErrorTitle = "错误"
$v3.ErrorMessage = "必须大于或者等于0"
$v3.ShowInput = $true
$v3.ShowError = $true

# P2:P1048576: whole numbers >= 0.
$v4 = $ws.Range("P2:P1048576").Validation
$v4.Add(1, 1, 7, 0)
$v4.ErrorTitle = "错误"
$v4.ErrorMessage = "必须大于或者等于0"
$v4.ShowInput = $true
$v4.ShowError = $true

# ---------------------------------------------------------------------------
# 6) Misc sheet cosmetics to match the authored workbook state.
# ---------------------------------------------------------------------------
# Drop the stale outline-level-row bookkeeping (no grouped rows remain).
$ws.Rows("1:1").OutlineLevel = 0

# Selection left where the author left it.
$ws.Range("K11").Select()
